$wb = $excel.ActiveWorkbook

# 1. Rename the "CarbonTrajectories" sheet to "CarbonTrajectory"
#    (fix sheet name to match the singular form used in parameters.jl code)
$carbonSheet = $wb.Worksheets.Item("CarbonTrajectories")
$carbonSheet.Name = "CarbonTrajectory"

# 2. Update the active selection/tab state on a few sheets to reflect where
#    the user ended up after making the edit.
[void]$carbonSheet.Range("C23").Select()

$energyTransitions = $wb.Worksheets.Item("EnergyTransitions")
[void]$energyTransitions.Activate()
$excel.ActiveWindow.Zoom = 156

$tempScenarios = $wb.Worksheets.Item("TemperatureScenarios")
[void]$tempScenarios.Activate()
[void]$tempScenarios.Range("C3").Select()
